# Transfer from Excel DB to JSON files for user details
# Append a new user record as row 4 to the Accounts sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "jsonTestv1"

# Password column holds numeric-looking text; force text so it is
# stored as a shared string rather than a number.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "12345"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").Value = "testone"
$ws.Range("D4").Value = "testone"

# DOB column is a date-formatted text string; force text so Excel
# does not reinterpret it as a date serial number.
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "09/18/2006"
$ws.Range("E4").Style = "Normal"

$ws.Range("F4").Value = "TT09182006"
$ws.Range("G4").Value = "jsonTestv1@gmail.com"
$ws.Range("H4").Value = "test st., Testv1, texas, US - 67676"

# Phone column holds a numeric-looking text string; force text.
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1234567890"
$ws.Range("I4").Style = "Normal"

$ws.Range("J4").Value = 69.0
